$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sub_materials_database")

$ws.Range("A1").Value = "sub_material_name"
$ws.Range("B1").Value = "chemical_composition"
